$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: now holds the new "RNN Model 4" entry (was "LSTM Model 1*") ---
# Gains the bottom-border separator formatting that used to sit under row 5.
$ws.Range("B6").Value = "RNN Model 4"
$ws.Range("C6").Value = "No"
$ws.Range("D6").Value = 0.6
$ws.Range("E6").Value = 0.67
$ws.Range("F6").Value = 0.41

# --- Row 7: now holds "LSTM Model 1*" (previously "MC-LSTM", removed from the table) ---
$ws.Range("B7").Value = "LSTM Model 1*"
$ws.Range("C7").Value = "Yes"
$ws.Range("D7").Value = 0.82
$ws.Range("E7").Value = 0.75
$ws.Range("F7").Value = 0.65

# Move the group-separator bottom border down from row 5 to row 6, so it now
# sits under the newly-added "RNN Model 4" row instead of "LSTM Model 3".
$ws.Range("B5:F5").Borders.Item(9).LineStyle = -4142

$ws.Range("B6:F6").Borders.Item(9).Color = 0
$ws.Range("B6:F6").Borders.Item(9).Weight = 2
$ws.Range("B6:F6").Borders.Item(9).LineStyle = 1
